$wb = $excel.ActiveWorkbook

# --- Locations sheet: remove the "AB" row ---
$locations = $wb.Worksheets.Item("Locations")
$locations.Rows.Item(3).Delete()

# --- AssetTypes sheet: update row 2's color, clear its link, and remove the AB/Cableway row ---
$assetTypes = $wb.Worksheets.Item("AssetTypes")
$assetTypes.Range("D2").Value = "#ccf610"
$assetTypes.Range("E2").ClearContents()
$assetTypes.Rows.Item(3).Delete()
